$wb = $excel.ActiveWorkbook
$wsTasks = $wb.Worksheets.Item("Tasks")
$wsPlanning = $wb.Worksheets.Item("Matt Planning")

# --- "Matt Planning" sheet: Project Backlog list maintenance ---
# Bug fix note: FieldAttrPanel would not load saved settings (nDec/FieldAttr)
# when the input Buffered Data Table changed.

# Remove the completed backlog item about reloading csv (row 14), which also
# removes its stray helper cell (C14) and its associated note (E14).
$wsPlanning.Rows.Item(14).Delete()

# Remove the completed backlog item "Default Table Name should be selected
# initially" (now at row 22 after the previous delete).
$wsPlanning.Rows.Item(22).Delete()

# Add new backlog items / notes for the qvx writer image work.
$wsPlanning.Cells.Item(22, 5).Value = "Advanced Settings Panel"
$wsPlanning.Cells.Item(23, 5).Value = 'Image of the word "qvx"'

# Update the remembered selections on each sheet.
$wsPlanning.Activate()
$wsPlanning.Range("A15").Select()

$wsTasks.Activate()
$wsTasks.Range("C16").Select()
